$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (current A,B become C,D)
$ws.Range("A:B").Insert()

# Header row
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

# Copy the style of the existing header cell (C1, which was A1) onto the new header cells
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 new values (existing product row, now shifted to C2:I2)
$ws.Range("A2").Value = "30/07/2024"
$ws.Range("B2").Value = "duartegabriella20230222153721"
$ws.Range("I2").Value = "https://produto.mercadolivre.com.br/MLB-4234174824-processador-jfa-digital-j4-redline-profissional-equalizador-_JM#position%3D1%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De13f5275-8729-4524-8243-4fb2341ca312"

# New row 3 with second product
$ws.Range("A3").Value = "30/07/2024"
$ws.Range("B3").Value = "duartegabriella20230222153721"
$ws.Range("C3").Value = "Processador De Audio Digital Equalizador Jfa J4 Red Line"
$ws.Range("D3").Value = "Sem Modelo"
$ws.Range("E3").Value = 399
# F3 mirrors the original sheet's blank "politica" cell (empty string, no value)
$ws.Range("G3").Value = "NA"
$ws.Range("H3").Value = "classico"
$ws.Range("I3").Value = "https://produto.mercadolivre.com.br/MLB-3226816467-processador-de-audio-digital-equalizador-jfa-j4-red-line-_JM#position%3D2%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De13f5275-8729-4524-8243-4fb2341ca312"

Write-Output ("UsedRange: " + $ws.UsedRange.Address())
Write-Output ("A1=" + $ws.Range("A1").Value2 + " B1=" + $ws.Range("B1").Value2 + " C1=" + $ws.Range("C1").Value2)
Write-Output ("A2=" + $ws.Range("A2").Value2 + " B2=" + $ws.Range("B2").Value2)
Write-Output ("I2=" + $ws.Range("I2").Value2)
Write-Output ("Row3: " + $ws.Range("A3").Value2 + " | " + $ws.Range("C3").Value2 + " | " + $ws.Range("I3").Value2)
